# Chose the best param for linear SVM
# - Sheet1 ("Тест по выборке из 5000"): the SVM(400) row is moved out to
#   Sheet2 ("Остальные"), the SVM(5000) row becomes the chosen leader
#   ("SVM (SVC) (лидер)") and every row after it shifts up by one, leaving
#   5 data rows instead of 6.
# - Sheet2 ("Остальные"): gains the former SVM(400) row and the previously
#   blank last row is filled in with the kNN(250000, 12) result.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# ---- Sheet1: "Тест по выборке из 5000" ----------------------------------

# Row 2 -> SVM (SVC) becomes the leader, now reporting the 5000-sample run.
$ws1.Range("A2").Value = "SVM (SVC) (лидер)"
$ws1.Range("B2").Value = "Linear, C=1"
$ws1.Range("C2").Value = 5000
$ws1.Range("D2").Value = "CV, 5"
$ws1.Range("E2").Value = 0.73
$ws1.Range("F2").Value = 0.02
$ws1.Range("G2").Value = 2622.81

# Row 3 -> kNN (was row 4)
$ws1.Range("A3").Value = "kNN"
$ws1.Range("B3").Value = 5
$ws1.Range("C3").Value = 5000
$ws1.Range("D3").Value = "CV, 5"
$ws1.Range("E3").Value = 0.76
$ws1.Range("F3").Value = 0.02
$ws1.Range("G3").Value = 0.49

# Row 4 -> kNN (лидер) (was row 5)
$ws1.Range("A4").Value = "kNN (лидер)"
$ws1.Range("B4").Value = 12
$ws1.Range("C4").Value = 5000
$ws1.Range("D4").Value = "CV, 5"
$ws1.Range("E4").Value = 0.78
$ws1.Range("F4").Value = 0.02
$ws1.Range("G4").Value = 0.65

# Row 5 -> RandomizedPCA+SVM (was row 6)
$ws1.Range("A5").Value = "RandomizedPCA+SVM"
$ws1.Range("B5").Value = "PCA(27, true), SVM(Linear, C=1)"
$ws1.Range("C5").Value = 5000
$ws1.Range("D5").Value = "CV, 5"
$ws1.Range("E5").Value = 0.74
$ws1.Range("F5").Value = 0.01
$ws1.Range("G5").Value = 6.87

# The old row 6 (previous SVM(400) row, now moved to Sheet2) is removed.
$ws1.Rows(6).Delete()

# ---- Sheet2: "Остальные" -------------------------------------------------

# Header row shrinks back down now that the SVM row below it is short.
$ws2.Rows(1).RowHeight = 30

# Row 2 -> SVM (SVC), 400 samples (moved from Sheet1 row 2)
$ws2.Range("A2").Value = "SVM (SVC)"
$ws2.Range("B2").Value = "Linear, C=1"
$ws2.Range("C2").Value = 400
$ws2.Range("D2").Value = "CV, 5"
$ws2.Range("E2").Value = 0.68
$ws2.Range("F2").Value = 0.06
$ws2.Range("G2").Value = 228.79

# Row 3 -> kNN, 250000 samples, k=5 (was row 2)
$ws2.Range("A3").Value = "kNN"
$ws2.Range("B3").Value = 5
$ws2.Range("C3").Value = 250000
$ws2.Range("D3").Value = "CV, 5"
$ws2.Range("E3").Value = 0.8
$ws2.Range("F3").Value = 0
$ws2.Range("G3").Value = 891.96

# Row 4 -> kNN, 250000 samples, k=12 (was row 3, previously-blank row 4 dropped)
$ws2.Range("A4").Value = "kNN"
$ws2.Range("B4").Value = 12
$ws2.Range("C4").Value = 250000
$ws2.Range("D4").Value = "CV, 5"
$ws2.Range("E4").Value = 0.81
$ws2.Range("F4").Value = 0
$ws2.Range("G4").Value = 920.64

# ---- Selections ----------------------------------------------------------
# Set Sheet2's selection first, then re-activate Sheet1 last so it remains
# the active / tab-selected sheet, matching the saved view state.
$ws2.Activate() | Out-Null
$ws2.Range("A2:A3").EntireRow.Select() | Out-Null

$ws1.Activate() | Out-Null
$ws1.Rows(2).Select() | Out-Null
